$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1.7188091074009018
$ws.Range("B2").Value = 0.36700743467069175
$ws.Range("C2").Value = -54.200260732023125
$ws.Range("D2").Value = 3.5524155020463346
$ws.Range("E2").Value = 509.78706515910977
$ws.Range("F2").Value = 8.6483402685236808
$ws.Range("G2").Value = 660.57699420414178
$ws.Range("H2").Value = 7.7851433262835137
$ws.Range("I2").Value = 7.8752251000398106
$ws.Range("J2").Value = 1.1250321571485444
$ws.Range("K2").Value = 0.45460950501875674
$ws.Range("L2").Value = 0.45460950501875674
$ws.Range("M2").Value = 0.63581382624014027
$ws.Range("N2").Value = -0.97228988306570696
